$wb = $excel.ActiveWorkbook

# --- 1. Unhide the "Options" sheet ---
$options = $wb.Worksheets.Item("Options")
$options.Visible = -1

# --- 2. Update the Options sheet's Category list (column C) ---
#    - "Post-event" removed
#    - "Health & Wellness" inserted after "Food & Drink"
#    - "Science & Education" inserted after "Sports & Recreation"
$categories = @(
    "Nature & Wildlife",
    "People & Lifestyle",
    "Business & Technology",
    "Travel & Places",
    "Food & Drink",
    "Health & Wellness",
    "Arts & Culture",
    "Sports & Recreation",
    "Science & Education",
    "Fashion & Beauty",
    "Interiors",
    "Abstract & Conceptual"
)
for ($i = 0; $i -lt $categories.Length; $i++) {
    $options.Cells.Item($i + 1, 3).Value = $categories[$i]
}

# --- 3. Extend the Images sheet's Category data validation to the new range ---
$images = $wb.Worksheets.Item("Images")
$catValidation = $images.Range("E2:E16")
$catValidation.Validation.Modify(3, 1, 1, "Options!`$C`$1:`$C`$12")
